$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two "Col/Row" labels in the "Worksheet (rows and columns)" /
# "Workbook (multiple worksheets)" mini-table (rows 29 and 31): split the
# combined "Col/Row" strings into a separate "Col" and "Row" label.
$ws.Range("D29").Value = "Col in a spreadsheet"
$ws.Range("D31").Value = "Row in a spreadsheet"

# Update the view: scroll the window so row 21 is back at the top and move
# the active selection to E31.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("E31").Select()
